$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - MIN
$ws.Range("C11:D11").Formula = "=MIN(C4:C8)"

# Row 12 - MAX
$ws.Range("C12:D12").Formula = "=MAX(C4:C8)"

# Row 13 - AVERAGE
$ws.Range("B13").Formula = "=AVERAGE(B4:B8)"
$ws.Range("C13:D13").Formula = "=AVERAGE(C4:C8)"

# Update the active selection to match the final state of the edit
$null = $ws.Range("I7").Select()
